# Basketball Tournament Rules.docx — apply yellow highlighting to a set of
# bullet paragraphs (and split one paragraph's run to insert new text).
#
# Strategy: Word's "highlight whole paragraph" (selecting from the very
# start of the paragraph through to -- and including -- its paragraph mark)
# updates both the run-level <w:rPr> AND the paragraph-mark's <w:pPr><w:rPr>,
# matching how Word itself records a highlighted bullet line. We therefore
# always grab Paragraphs(n).Range (which already spans start..end inclusive
# of the pilcrow) and set .Font.HighlightColorIndex = wdYellow (7) on it.

$d = $word.ActiveDocument
$wdYellow = 7

function Highlight-Paragraph($index) {
    $p = $d.Paragraphs($index)
    $p.Range.Font.HighlightColorIndex = $wdYellow
}

# --- "Every win/loss must be recorded after each individual game" ---------
# This paragraph also gains new text ("and score ") in the middle, splitting
# the original single run into three. Insert the text first (while using a
# throwaway Bold toggle to force Word to keep the insertion as its own run
# instead of silently re-merging it with identical neighbouring runs), and
# only then apply the highlight across the whole paragraph so every run
# (old and new) -- and the paragraph mark -- end up highlighted.
$find = $d.Content
$found = $find.Find.Execute("Every win/loss ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Every win/loss ' to split/insert into"
}
$insertPoint = $find.End
$insertRange = $d.Range($insertPoint, $insertPoint)
$insertRange.InsertAfter("and score ")

$newRunRange = $d.Range($insertPoint, $insertPoint + 10)
$newRunRange.Bold = 1
$newRunRange.Bold = 0

Highlight-Paragraph 9

# --- Simple whole-paragraph highlights -------------------------------------
$simpleParagraphIndexes = @(5, 6, 10, 11, 21, 22, 32, 37, 40, 41, 42, 43, 44)
foreach ($idx in $simpleParagraphIndexes) {
    Highlight-Paragraph $idx
}
